$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "I ordered a few items in Order ID 1263. Have they been shipped yet?"
$ws.Range("C2").Value = "{'orders': [{'Product ID': 43491, 'ProductName': 'bosch kir41af30g 122cm integrated fridge white', 'Category': 'Fridges', 'Category ID': 2623, 'OrderID': 1263, 'CustomerID': 726, 'OrderStatus': 'Shipped', 'ReturnEligible': False, 'ShippingDate': '2024-10-18 12:24:37.765490'}]}"
$ws.Range("B3").Value = "I placed Order ID 3179 a while ago. Could you please check on the shipping status of each item?"
$ws.Range("C3").Value = "{'orders': [{'Product ID': 46508, 'ProductName': 'liebherr k hlschrankl premium mit biofresh kbpes 4354', 'Category': 'Fridges', 'Category ID': 2623, 'OrderID': 3179, 'CustomerID': 2161, 'OrderStatus': 'Cancelled', 'ReturnEligible': False, 'ShippingDate': '2024-09-29 12:24:37.765490'}]}"
$ws.Range("B4").Value = "Order ID 2193 contains several items. Can you tell me if each one has shipped?"
$ws.Range("C4").Value = "{'orders': [{'Product ID': 16185, 'ProductName': 'hp intel xeon e7530 1.866ghz 12mb l3 processor', 'Category': 'CPUs', 'Category ID': 2615, 'OrderID': 2193, 'CustomerID': 2408, 'OrderStatus': 'Pending', 'ReturnEligible': False, 'ShippingDate': '2024-10-26 12:24:37.765490'}]}"
$ws.Range("B5").Value = "In Order ID 1182, I purchased multiple items. I'd like to know the shipping status for each of them."
$ws.Range("C5").Value = "{'orders': [{'Product ID': 44483, 'ProductName': 'zanussi zba15021sv integrated fridge', 'Category': 'Fridges', 'Category ID': 2623, 'OrderID': 1182, 'CustomerID': 193, 'OrderStatus': 'Pending', 'ReturnEligible': False, 'ShippingDate': '2024-10-05 12:24:37.765490'}]}"
$ws.Range("B6").Value = "I have multiple items under Order ID 1384. Can you check the shipping status for each item, please?"
$ws.Range("C6").Value = "{'orders': [{'Product ID': 42483, 'ProductName': 'siemens k hl gefrier k hl kombination wei kg39vvw31 eek a', 'Category': 'Fridge Freezers', 'Category ID': 2622, 'OrderID': 1384, 'CustomerID': 2059, 'OrderStatus': 'Pending', 'ReturnEligible': False, 'ShippingDate': '2024-10-08 12:24:37.765490'}]}"
$ws.Range("B7").Value = "I ordered multiple items in Order ID 1530. Can you check if each item has shipped?"
$ws.Range("C7").Value = "{'orders': [{'Product ID': 46062, 'ProductName': 'bosch k hlger t kil82ad40', 'Category': 'Fridges', 'Category ID': 2623, 'OrderID': 2830, 'CustomerID': 1349, 'OrderStatus': 'Shipped', 'ReturnEligible': False, 'ShippingDate': '2024-10-09 12:24:37.765490'}]}"
$ws.Range("B8").Value = "I placed an order with multiple items (Order ID 3183). Could you please tell me the shipping status of each item?"
$ws.Range("C8").Value = "{'orders': [{'Product ID': 16829, 'ProductName': 'intel xeon e7340 processor 2.4 ghz box 8 mb l2', 'Category': 'CPUs', 'Category ID': 2615, 'OrderID': 2946, 'CustomerID': 119, 'OrderStatus': 'Shipped', 'ReturnEligible': False, 'ShippingDate': '2024-10-13 12:24:37.765490'}]}"
$ws.Range("B9").Value = "Order 1113 has several items. Have they all shipped yet?"
$ws.Range("C9").Value = "{'orders': [{'Product ID': 34428, 'ProductName': 'bosch waw32540', 'Category': 'Washing Machines', 'Category ID': 2620, 'OrderID': 2728, 'CustomerID': 294, 'OrderStatus': 'Cancelled', 'ReturnEligible': False, 'ShippingDate': '2024-10-13 12:24:37.765490'}]}"
$ws.Range("B10").Value = "I'd like to know the shipping status for each item in Order ID 1842. It contains multiple items."
$ws.Range("C10").Value = "{'orders': [{'Product ID': 16843, 'ProductName': 'hewlett packard enterprise amd opteron 2352 2.1ghz 2mb l3 processor', 'Category': 'CPUs', 'Category ID': 2615, 'OrderID': 1753, 'CustomerID': 1584, 'OrderStatus': 'Shipped', 'ReturnEligible': False, 'ShippingDate': '2024-10-25 12:24:37.765490'}]}"
$ws.Range("B11").Value = "Can you tell me if everything in Order 1948 has shipped? There are multiple items on this order."
$ws.Range("C11").Value = "{'orders': [{'Product ID': 43969, 'ProductName': 'bosch kil22vs30g serie 4 built in single door fridges 88cm height with ice box white', 'Category': 'Fridges', 'Category ID': 2623, 'OrderID': 1087, 'CustomerID': 1040, 'OrderStatus': 'Cancelled', 'ReturnEligible': False, 'ShippingDate': '2024-10-25 12:24:37.765490'}]}"
